$wb = $excel.ActiveWorkbook

# Rename the "Device" sheet to "Apparatus"
$ws = $wb.Worksheets.Item("Device")
$ws.Name = "Apparatus"

# Update the text on the renamed sheet: "Device" -> "Apparatus"
$ws.Range("B2").Value = "Apparatus type"
$ws.Range("C2").Value = "Apparatus parameters"
$ws.Range("A1").Value = "This sheet summarizes the apparatuses connected to buses."

# Move the active tab/selection from "Bus" to "Apparatus"
$ws.Activate()
$ws.Range("A2").Select()

$wb.Save()
